# Refresh the "paises" (COVID per-country) worksheet with updated figures.
#
# The underlying data is ranked by total cases, so when a country's
# numbers are refreshed its rank can change, which means the country name
# in column A moves to a different row together with its own stats while
# the countries it overtakes shift down by one row. Rows that are not
# impacted by any re-ranking only get their statistics columns refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Septiembre de 2020 a las 12:56"

# Rows where the country ranking changed (country name + stats move together)

# Emiratos Arabes Unidos overtakes Paises Bajos
$ws.Range("A44").Value = "Emiratos Arabes Unidos"
$ws.Range("B44").Value = 73471
$ws.Range("C44").Value = 705
$ws.Range("D44").Value = 63652
$ws.Range("E44").Value = 9431
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 388

$ws.Range("A45").Value = "Paises Bajos"  # takes over the old Emiratos Arabes Unidos row
$ws.Range("B45").Value = 73208
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 6237

# Nepal overtakes Ghana, Armenia and Kirguistan, pushing them down one row each
$ws.Range("A59").Value = "Nepal"
$ws.Range("B59").Value = 45277
$ws.Range("C59").Value = 1041
$ws.Range("D59").Value = 27127
$ws.Range("E59").Value = 17870
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 9
$ws.Range("H59").Value = 280

$ws.Range("A60").Value = "Ghana"
$ws.Range("B60").Value = 44777
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 43693
$ws.Range("E60").Value = 801
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 283

$ws.Range("A61").Value = "Armenia"
$ws.Range("B61").Value = 44649
$ws.Range("C61").Value = 188
$ws.Range("D61").Value = 39823
$ws.Range("E61").Value = 3931
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 895

$ws.Range("A62").Value = "Kirguistan"
$ws.Range("B62").Value = 44293
$ws.Range("C62").Value = 94
$ws.Range("D62").Value = 39599
$ws.Range("E62").Value = 3634
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 1060

# Malta overtakes Botsuana and Reunion, pushing them down one row each
$ws.Range("A146").Value = "Malta"
$ws.Range("B146").Value = 2014
$ws.Range("C146").Value = 30
$ws.Range("D146").Value = 1601
$ws.Range("E146").Value = 399
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 1
$ws.Range("H146").Value = 14

$ws.Range("A147").Value = "Botsuana"
$ws.Range("B147").Value = 2002
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 493
$ws.Range("E147").Value = 1501
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 8

$ws.Range("A148").Value = "Reunion"
$ws.Range("B148").Value = 2002
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 880
$ws.Range("E148").Value = 1112
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 10

# Gibraltar overtakes Mongolia, pushing it down one row
$ws.Range("A183").Value = "Gibraltar"
$ws.Range("B183").Value = 312
$ws.Range("C183").Value = 7
$ws.Range("D183").Value = 255
$ws.Range("E183").Value = 57
$ws.Range("F183").Value = 0

$ws.Range("A184").Value = "Mongolia"
$ws.Range("B184").Value = 310
$ws.Range("C184").Value = 0
$ws.Range("D184").Value = 296
$ws.Range("E184").Value = 14
$ws.Range("F184").Value = 0

# Rows where only statistics were refreshed (country unchanged)
$ws.Range("B17").Value = 323565
$ws.Range("C17").Value = 1950
$ws.Range("D17").Value = 217852
$ws.Range("E17").Value = 101266
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = 4447

$ws.Range("B30").Value = 119864
$ws.Range("C30").Value = 227
$ws.Range("D30").Value = 116780
$ws.Range("E30").Value = 2882
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 202

$ws.Range("B37").Value = 93864
$ws.Range("C37").Value = 1269
$ws.Range("D37").Value = 40028
$ws.Range("E37").Value = 49986
$ws.Range("G37").Value = 38
$ws.Range("H37").Value = 3850

$ws.Range("E55").Value = 3576
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 196

$ws.Range("B63").Value = 43957
$ws.Range("C63").Value = 425
$ws.Range("E63").Value = 5444

$ws.Range("B75").Value = 25575
$ws.Range("C75").Value = 433
$ws.Range("D75").Value = 16661
$ws.Range("E75").Value = 8738
$ws.Range("G75").Value = 6
$ws.Range("H75").Value = 176

$ws.Range("B87").Value = 13948
$ws.Range("C87").Value = 67
$ws.Range("D87").Value = 9851
$ws.Range("E87").Value = 3807
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 290

$ws.Range("B112").Value = 4858
$ws.Range("C112").Value = 7
$ws.Range("D112").Value = 4492
$ws.Range("E112").Value = 272
